$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new feedback row (row 5) with the latest entry.
$ws.Range("A5").Value = "r456"
$ws.Range("B5").Value = "george"
$ws.Range("C5").Value = "tiny tim must go"
$ws.Range("D5").Value = "2025-09-30 20:26:40"
